$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1322").Value = "U45_01"
$ws.Range("B1322").Value = 45
$ws.Range("C1322").Value = "Cảm nhận, cảm giác"
$ws.Range("D1322").Value = "Sense"
$ws.Range("E1322").Value = "Success gave me a sense of happiness"
$ws.Range("F1322").Value = "a sense of something"
$ws.Range("G1322").Value = "N"

$ws.Range("A1323").Value = "U45_02"
$ws.Range("B1323").Value = 45
$ws.Range("C1323").Value = "Quan điểm"
$ws.Range("D1323").Value = "Opinion"
$ws.Range("E1323").Value = "He has opinions about everything"
$ws.Range("F1323").Value = "have an opinion about something"
$ws.Range("G1323").Value = "N"

$ws.Range("A1324").Value = "U45_03"
$ws.Range("B1324").Value = 45
$ws.Range("C1324").Value = "Hữu ích"
$ws.Range("D1324").Value = "Useful"
$ws.Range("E1324").Value = "Your insight (hiểu biết của cậu) was useful for my research"
$ws.Range("F1324").Value = "useful for somebody or something"
$ws.Range("G1324").Value = "Adj"

$ws.Range("A1325").Value = "U45_04"
$ws.Range("B1325").Value = 45
$ws.Range("C1325").Value = "Có lẽ"
$ws.Range("D1325").Value = "Perhaps"
$ws.Range("E1325").Value = "Perhaps we should tell him the truth"
$ws.Range("F1325").Value = "Perhaps somebody should do something"
$ws.Range("G1325").Value = "Adv"

$ws.Range("A1326").Value = "U45_05"
$ws.Range("B1326").Value = 45
$ws.Range("C1326").Value = "Chắc chắn"
$ws.Range("D1326").Value = "Definitely"
$ws.Range("E1326").Value = "This is definitely not her first time lying."
$ws.Range("F1326").Value = "definitely not / chắc chắn không phải"
$ws.Range("G1326").Value = "Adv"

$ws.Range("A1327").Value = "U45_06"
$ws.Range("B1327").Value = 45
$ws.Range("C1327").Value = "Sự thật"
$ws.Range("D1327").Value = "Truth"
$ws.Range("E1327").Value = "The truth is, I lied to you"
$ws.Range("F1327").Value = "the truth is / sự thật là"
$ws.Range("G1327").Value = "N"

$ws.Range("A1328").Value = "U45_07"
$ws.Range("B1328").Value = 45
$ws.Range("C1328").Value = "Bình luận"
$ws.Range("D1328").Value = "Comment"
$ws.Range("E1328").Value = "Please leave a comment if you have any questions"
$ws.Range("F1328").Value = "leave a comment / để lại bình luận"
$ws.Range("G1328").Value = "N"

$ws.Range("A1329").Value = "U45_08"
$ws.Range("B1329").Value = 45
$ws.Range("C1329").Value = "Tin"
$ws.Range("D1329").Value = "Believe"
$ws.Range("E1329").Value = "Do you believe in me"
$ws.Range("F1329").Value = "believe in somebody / tin vào ai đó"
$ws.Range("G1329").Value = "V"

$ws.Range("A1330").Value = "U45_09"
$ws.Range("B1330").Value = 45
$ws.Range("C1330").Value = "Hoàn toàn"
$ws.Range("D1330").Value = "Totally"
$ws.Range("E1330").Value = "What you said was totally unacceptable"
$ws.Range("F1330").Value = "totally unacceptable / hoàn toàn không thể chấp nhận được"
$ws.Range("G1330").Value = "Adv"

$ws.Range("A1331").Value = "U45_10"
$ws.Range("B1331").Value = 45
$ws.Range("C1331").Value = "Thấy phiền"
$ws.Range("D1331").Value = "Mind"
$ws.Range("E1331").Value = "Would you mind if I ask you something?"
$ws.Range("F1331").Value = "would you mind / bạn có phiền"
$ws.Range("G1331").Value = "V"

$ws.Range("A1332").Value = "U45_11"
$ws.Range("B1332").Value = 45
$ws.Range("C1332").Value = "Sự tự do"
$ws.Range("D1332").Value = "Freedom"
$ws.Range("E1332").Value = "Children need the freedom to make mistakes."
$ws.Range("F1332").Value = "freedom to do something "
$ws.Range("G1332").Value = "N"

$ws.Range("A1333").Value = "U45_12"
$ws.Range("B1333").Value = 45
$ws.Range("C1333").Value = "Lý lẽ, cuộc tranh cãi"
$ws.Range("D1333").Value = "Argument"
$ws.Range("E1333").Value = "I had an argument with her yesterday"
$ws.Range("F1333").Value = "an argument with somebody / một cuộc tranh luận với ai đó"
$ws.Range("G1333").Value = "N"

$ws.Range("A1334").Value = "U45_13"
$ws.Range("B1334").Value = 45
$ws.Range("C1334").Value = "Cuộc tranh biện"
$ws.Range("D1334").Value = "debate"
$ws.Range("E1334").Value = "We were having a debate on justice (công lý)"
$ws.Range("F1334").Value = "a debate on something"
$ws.Range("G1334").Value = "N"

$ws.Range("A1335").Value = "U45_14"
$ws.Range("B1335").Value = 45
$ws.Range("C1335").Value = "Dường như"
$ws.Range("D1335").Value = "Seem"
$ws.Range("E1335").Value = "They seem to be scared at first (lúc đầu họ dường như là bị sợ sệt)"
$ws.Range("F1335").Value = "seem to be / dường như là"
$ws.Range("G1335").Value = "V"

$ws.Range("A1336").Value = "U45_15"
$ws.Range("B1336").Value = 45
$ws.Range("C1336").Value = "Sự lựa chọn"
$ws.Range("D1336").Value = "Choice"
$ws.Range("E1336").Value = "Make a choice between apples and oranges"
$ws.Range("F1336").Value = "choice between A and B"
$ws.Range("G1336").Value = "N"

$ws.Range("A1337").Value = "U45_16"
$ws.Range("B1337").Value = 45
$ws.Range("C1337").Value = "Đồng ý"
$ws.Range("D1337").Value = "Agree"
$ws.Range("E1337").Value = "I agree with her on moving forward (tiếp tục triển khai)"
$ws.Range("F1337").Value = "agree with somebody"
$ws.Range("G1337").Value = "V"

$ws.Range("A1338").Value = "U45_17"
$ws.Range("B1338").Value = 45
$ws.Range("C1338").Value = "Thực tế"
$ws.Range("D1338").Value = "Realistic"
$ws.Range("E1338").Value = "Be realistic about your expectations / Hãy thực tế về những sự kỳ vọng của mình"
$ws.Range("F1338").Value = "To be realistic about something"
$ws.Range("G1338").Value = "Adj"

$ws.Range("A1339").Value = "U45_18"
$ws.Range("B1339").Value = 45
$ws.Range("C1339").Value = "Chia"
$ws.Range("D1339").Value = "Divide"
$ws.Range("E1339").Value = "Divide the apples into 4 servings (phần)"
$ws.Range("F1339").Value = "divide into something / chia thành cái gì đó"
$ws.Range("G1339").Value = "V"

$ws.Range("A1340").Value = "U45_19"
$ws.Range("B1340").Value = 45
$ws.Range("C1340").Value = "Công bằng"
$ws.Range("D1340").Value = "fair"
$ws.Range("E1340").Value = "My teacher is fair to all students"
$ws.Range("F1340").Value = "to be fair to somebody / công bằng với ai đó"
$ws.Range("G1340").Value = "Adj"

$ws.Range("A1341").Value = "U45_20"
$ws.Range("B1341").Value = 45
$ws.Range("C1341").Value = "Nhân chứng"
$ws.Range("D1341").Value = "Witness"
$ws.Range("E1341").Value = "I was a witness to his crime"
$ws.Range("F1341").Value = "a witness to something"
$ws.Range("G1341").Value = "N"

$ws.Range("A1342").Value = "U45_21"
$ws.Range("B1342").Value = 45
$ws.Range("C1342").Value = "Phiên bản"
$ws.Range("D1342").Value = "Version"
$ws.Range("E1342").Value = "Here is the latest version of the iPhone."
$ws.Range("F1342").Value = "version of something"
$ws.Range("G1342").Value = "N"

$ws.Range("A1343").Value = "U45_22"
$ws.Range("B1343").Value = 45
$ws.Range("C1343").Value = "Lời tuyên bố"
$ws.Range("D1343").Value = "statement"
$ws.Range("E1343").Value = "The film was a statement about fairness."
$ws.Range("F1343").Value = "statement about something"
$ws.Range("G1343").Value = "N"

$ws.Range("A1344").Value = "U45_23"
$ws.Range("B1344").Value = 45
$ws.Range("C1344").Value = "Về mặt cá nhân"
$ws.Range("D1344").Value = "Personally"
$ws.Range("E1344").Value = "Personally speaking, I feel offended by the question (câu hỏi có phần xúc phạm)"
$ws.Range("F1344").Value = "personally speaking / về mặt cá nhân mà nói"
$ws.Range("G1344").Value = "Adv"

$ws.Range("A1345").Value = "U45_24"
$ws.Range("B1345").Value = 45
$ws.Range("C1345").Value = "Quan điểm"
$ws.Range("D1345").Value = "View"
$ws.Range("E1345").Value = "From my point of view, you seem nice"
$ws.Range("F1345").Value = "from my point of view / theo quan điểm của tôi"
$ws.Range("G1345").Value = "N"

$ws.Range("A1346").Value = "U45_25"
$ws.Range("B1346").Value = 45
$ws.Range("C1346").Value = "Lý thuyết"
$ws.Range("D1346").Value = "Theory"
$ws.Range("E1346").Value = "Hawking discovered the theory of black holes"
$ws.Range("F1346").Value = "Theory of something"
$ws.Range("G1346").Value = "N"

$ws.Range("A1347").Value = "U45_26"
$ws.Range("B1347").Value = 45
$ws.Range("C1347").Value = "Cho là, tự nhận là"
$ws.Range("D1347").Value = "Claim"
$ws.Range("E1347").Value = "He claims to be an honest man"
$ws.Range("F1347").Value = "claim to be something"
$ws.Range("G1347").Value = "V"

$ws.Range("A1348").Value = "U45_27"
$ws.Range("B1348").Value = 45
$ws.Range("C1348").Value = "Lạc quan"
$ws.Range("D1348").Value = "Optimistic"
$ws.Range("E1348").Value = "We are optimistic about the future"
$ws.Range("F1348").Value = "optimistic about something"
$ws.Range("G1348").Value = "Adj"

$ws.Range("A1349").Value = "U45_28"
$ws.Range("B1349").Value = 45
$ws.Range("C1349").Value = "Khá"
$ws.Range("D1349").Value = "Quite"
$ws.Range("E1349").Value = "The meal was quired good"
$ws.Range("F1349").Value = "quite good/ khá tốt, khá ngon"
$ws.Range("G1349").Value = "Adv"

$ws.Range("A1350").Value = "U45_29"
$ws.Range("B1350").Value = 45
$ws.Range("C1350").Value = "Đáng kinh ngạc, khó tin"
$ws.Range("D1350").Value = "Incredible"
$ws.Range("E1350").Value = "He told us an incredible story"
$ws.Range("F1350").Value = "an incredible story"
$ws.Range("G1350").Value = "Adj"

$ws.Range("A1351").Value = "U45_30"
$ws.Range("B1351").Value = 45
$ws.Range("C1351").Value = "Sự đồng cảm"
$ws.Range("D1351").Value = "Sympathy"
$ws.Range("E1351").Value = "Have sympathy for her loss (sự mất mát của cô ấy)"
$ws.Range("F1351").Value = "Have sympathy for somebody"
$ws.Range("G1351").Value = "N"

$ws.Range("E1350").Select()
